# This script applies a set of stock-report corrections to the single
# worksheet in the workbook. Some item rows had their order-quantity
# (column F) corrected, which also changes the corresponding line value
# (column G = rate * quantity). A handful of adjacent item rows (same
# product, different batch/rate) had their batch code (B), rate (E),
# quantity (F) and value (G) swapped between the two rows. All of the
# "Sub Total:" rows (column B) and the final "Sub Total:"/"Grand Total:"
# rows at the bottom of the sheet are simple sums of the values above
# them, and are updated here to stay consistent with the corrected
# quantities.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Item rows: corrected quantity (F) and recalculated value (G) ---
$ws.Range("F71").Value = 312
$ws.Range("G71").Value = 19874.4

$ws.Range("F115").Value = 187
$ws.Range("G115").Value = 18103.47

# Row 219 / 220: batch (B), rate (E), quantity (F) and value (G) swapped
$ws.Range("B219").Value = 61610
$ws.Range("E219").Value = 122.71
$ws.Range("F219").Value = -58
$ws.Range("G219").Value = -5957.18

$ws.Range("B220").Value = 63565
$ws.Range("E220").Value = 109.19
$ws.Range("F220").Value = 60
$ws.Range("G220").Value = 6162.6

# Row 232 / 233: batch (B), rate (E), quantity (F) and value (G) swapped
$ws.Range("B232").Value = 63510
$ws.Range("E232").Value = 50.66
$ws.Range("F232").Value = 114
$ws.Range("G232").Value = 5430.96

$ws.Range("B233").Value = 55356
$ws.Range("E233").Value = 54.04
$ws.Range("F233").Value = -158
$ws.Range("G233").Value = -7527.12

$ws.Range("F234").Value = 36
$ws.Range("G234").Value = 1847.52

# Row 243 / 244: batch (B), rate (E), quantity (F) and value (G) swapped
$ws.Range("B243").Value = 63560
$ws.Range("E243").Value = 134.87
$ws.Range("F243").Value = 1
$ws.Range("G243").Value = 126.86

$ws.Range("B244").Value = 60325
$ws.Range("E244").Value = 151.57
$ws.Range("F244").Value = -102
$ws.Range("G244").Value = -12939.72

$ws.Range("F249").Value = 137
$ws.Range("G249").Value = 18881.34

$ws.Range("F278").Value = 9
$ws.Range("G278").Value = 1235.88

$ws.Range("F296").Value = 42
$ws.Range("G296").Value = 890.4

$ws.Range("F300").Value = 162
$ws.Range("G300").Value = 20235.42

$ws.Range("F328").Value = 37
$ws.Range("G328").Value = 1376.77

$ws.Range("F354").Value = 11
$ws.Range("G354").Value = 754.49

# Row 385 / 386: batch (B), rate (E), quantity (F) and value (G) swapped
$ws.Range("B385").Value = 65067
$ws.Range("E385").Value = 15.65
$ws.Range("F385").Value = 126
$ws.Range("G385").Value = 1855.98

$ws.Range("B386").Value = 53595
$ws.Range("E386").Value = 17.61
$ws.Range("F386").Value = -335
$ws.Range("G386").Value = -4934.55

# Row 442 / 443: batch (B), rate (E), quantity (F) and value (G) swapped
$ws.Range("B442").Value = 53319
$ws.Range("E442").Value = 310.64
$ws.Range("F442").Value = -6
$ws.Range("G442").Value = -1643.52

$ws.Range("B443").Value = 64810
$ws.Range("E443").Value = 291.22
$ws.Range("F443").Value = 4
$ws.Range("G443").Value = 1095.68

$ws.Range("F453").Value = 17
$ws.Range("G453").Value = 450.67

# Row 473 / 474: batch (B), rate (E), quantity (F) and value (G) swapped
$ws.Range("B473").Value = 60022
$ws.Range("E473").Value = 37.22
$ws.Range("F473").Value = -113
$ws.Range("G473").Value = -3709.79

$ws.Range("B474").Value = 64830
$ws.Range("E474").Value = 34.9
$ws.Range("F474").Value = 107
$ws.Range("G474").Value = 3512.81

$ws.Range("F509").Value = 200
$ws.Range("G509").Value = 16076

$ws.Range("F573").Value = 18
$ws.Range("G573").Value = 735.66

$ws.Range("F578").Value = 63
$ws.Range("G578").Value = 3143.07

$ws.Range("F595").Value = 8
$ws.Range("G595").Value = 309.52

$ws.Range("F599").Value = 1379
$ws.Range("G599").Value = 224928.69

$ws.Range("F601").Value = 370
$ws.Range("G601").Value = 104661.9

$ws.Range("F602").Value = 316
$ws.Range("G602").Value = 45709.4

$ws.Range("F615").Value = 81
$ws.Range("G615").Value = 10153.35

# --- Sub Total / Grand Total rows: resulting sums recalculated ---
$ws.Range("B90").Value = 170046.48
$ws.Range("B117").Value = 11944.53
$ws.Range("B260").Value = 173031.33
$ws.Range("B304").Value = 166754.17
$ws.Range("B330").Value = 26711.11
$ws.Range("B358").Value = 34777.48
$ws.Range("B460").Value = 12379.77
$ws.Range("B510").Value = 21480.88
$ws.Range("B583").Value = 13375.44
$ws.Range("B596").Value = 309.52
$ws.Range("B606").Value = 376148.04
$ws.Range("B618").Value = 41990.3
$ws.Range("B619").Value = 1607315.9
$ws.Range("B620").Value = 1607315.9
